$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.04068733621174444, 1.347492170625331, 10.12087998327202, 3.181333051296583, 3.212726103237388, 51)
    3  = @(0.02310952642605121, 1.506637317013713, 9.664068608552252, 3.108708511351981, 3.140183043920715, 50)
    4  = @(0.04953615400542137, 1.37882258673296, 6.90966373676957, 2.628623924560067, 2.65539264753735, 49)
    5  = @(0.06730071474157288, 1.432547271739592, 8.491317398801943, 2.913986513146885, 2.944037698520868, 48)
    6  = @(0.03756714743281693, 1.591018637028787, 10.0026066766869, 3.162689785085931, 3.196656484066325, 47)
    7  = @(0.05601094792970936, 1.592722824606551, 8.729723085726032, 2.954610479526198, 2.986722281975306, 46)
    8  = @(0.0002899658432591508, 1.568355227381864, 8.737285131131101, 2.955889905109983, 2.989290836839259, 45)
    9  = @(-0.02005145808999941, 1.594452076336724, 9.27842226268784, 3.046050272514858, 3.081199135626279, 44)
    10 = @(0.01513648954317695, 1.65609166469574, 9.395576928003257, 3.065220534970242, 3.101458781312453, 43)
    11 = @(0.004535327475175432, 1.624160030100705, 9.43613802475137, 3.071829751915195, 3.109062022827647, 42)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
